$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2125
$ws.Range("C2").Value = 0.51875
$ws.Range("J2").Value = 0.0125
$ws.Range("P2").Value = 0.159375
$ws.Range("S2").Value = 0.096875
$ws.Range("B3").Value = 0.01169590643274854
$ws.Range("C3").Value = 0.02339181286549707
$ws.Range("J3").Value = 0.02339181286549707
$ws.Range("P3").Value = 0.7719298245614035
$ws.Range("S3").Value = 0.1695906432748538
$ws.Range("P4").Value = 0.8431372549019608
$ws.Range("S4").Value = 0.1568627450980392
$ws.Range("B6").Value = 0.03389830508474576
$ws.Range("D6").Value = 0.03813559322033899
$ws.Range("F6").Value = 0.05084745762711865
$ws.Range("J6").Value = 0.1864406779661017
$ws.Range("O6").Value = 0.008474576271186441
$ws.Range("Q6").Value = 0.1779661016949153
$ws.Range("R6").Value = 0.09745762711864407
$ws.Range("S6").Value = 0.4067796610169492
$ws.Range("B7").Value = 0.09677419354838709
$ws.Range("D7").Value = 0.02150537634408602
$ws.Range("E7").Value = 0.005376344086021506
$ws.Range("F7").Value = 0.04838709677419355
$ws.Range("J7").Value = 0.06989247311827956
$ws.Range("O7").Value = 0.03763440860215054
$ws.Range("Q7").Value = 0.1935483870967742
$ws.Range("R7").Value = 0.1129032258064516
$ws.Range("S7").Value = 0.4139784946236559
$ws.Range("B8").Value = 0.08228980322003578
$ws.Range("D8").Value = 0.01610017889087657
$ws.Range("E8").Value = 0.001788908765652952
$ws.Range("F8").Value = 0.07513416815742398
$ws.Range("J8").Value = 0.08586762075134168
$ws.Range("O8").Value = 0.01252236135957066
$ws.Range("Q8").Value = 0.184257602862254
$ws.Range("R8").Value = 0.10912343470483
$ws.Range("S8").Value = 0.4329159212880143
$ws.Range("B9").Value = 0.1129032258064516
$ws.Range("D9").Value = 0.02419354838709677
$ws.Range("F9").Value = 0.06854838709677419
$ws.Range("J9").Value = 0.08870967741935484
$ws.Range("O9").Value = 0.01209677419354839
$ws.Range("Q9").Value = 0.1491935483870968
$ws.Range("R9").Value = 0.1088709677419355
$ws.Range("S9").Value = 0.4354838709677419
$ws.Range("B10").Value = 0.1042415528396837
$ws.Range("D10").Value = 0.01941049604601006
$ws.Range("E10").Value = 0.0007189072609633358
$ws.Range("F10").Value = 0.07117181883537024
$ws.Range("J10").Value = 0.08123652048885693
$ws.Range("O10").Value = 0.01725377426312006
$ws.Range("Q10").Value = 0.2207045291157441
$ws.Range("R10").Value = 0.102803738317757
$ws.Range("S10").Value = 0.3824586628324946
$ws.Range("G11").Value = 0.1182795698924731
$ws.Range("J11").Value = 0.09677419354838709
$ws.Range("K11").Value = 0.1899641577060932
$ws.Range("L11").Value = 0.5770609318996416
$ws.Range("S11").Value = 0.01792114695340502
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.1785714285714286
$ws.Range("L12").Value = 0.02976190476190476
$ws.Range("S12").Value = 0.04166666666666666
$ws.Range("F15").Value = 0.01652892561983471
$ws.Range("H15").Value = 0.1694214876033058
$ws.Range("I15").Value = 0.05785123966942149
$ws.Range("J15").Value = 0.3553719008264463
$ws.Range("K15").Value = 0.04132231404958678
$ws.Range("M15").Value = 0.02066115702479339
$ws.Range("O15").Value = 0.0743801652892562
$ws.Range("S15").Value = 0.2644628099173554
$ws.Range("F16").Value = 0.0091324200913242
$ws.Range("H16").Value = 0.1735159817351598
$ws.Range("I16").Value = 0.0867579908675799
$ws.Range("J16").Value = 0.4337899543378995
$ws.Range("K16").Value = 0.091324200913242
$ws.Range("M16").Value = 0.0228310502283105
$ws.Range("O16").Value = 0.0776255707762557
$ws.Range("S16").Value = 0.1050228310502283
$ws.Range("F17").Value = 0.01535508637236084
$ws.Range("H17").Value = 0.2168905950095969
$ws.Range("I17").Value = 0.1151631477927063
$ws.Range("J17").Value = 0.3992322456813819
$ws.Range("K17").Value = 0.08253358925143954
$ws.Range("M17").Value = 0.007677543186180422
$ws.Range("N17").Value = 0.001919385796545105
$ws.Range("O17").Value = 0.05374280230326296
$ws.Range("S17").Value = 0.1074856046065259
$ws.Range("F18").Value = 0.01107011070110701
$ws.Range("H18").Value = 0.1771217712177122
$ws.Range("I18").Value = 0.1180811808118081
$ws.Range("J18").Value = 0.4280442804428044
$ws.Range("K18").Value = 0.07380073800738007
$ws.Range("M18").Value = 0.01476014760147601
$ws.Range("O18").Value = 0.07749077490774908
$ws.Range("S18").Value = 0.0996309963099631
$ws.Range("F19").Value = 0.007670850767085077
$ws.Range("H19").Value = 0.2273361227336123
$ws.Range("I19").Value = 0.08647140864714087
$ws.Range("J19").Value = 0.403068340306834
$ws.Range("K19").Value = 0.09274755927475593
$ws.Range("M19").Value = 0.01743375174337517
$ws.Range("N19").Value = 0.002789400278940028
$ws.Range("O19").Value = 0.05857740585774059
$ws.Range("S19").Value = 0.103905160390516
